$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cycle_2021-2022")

# --- Append the new "44555" data block (rows 89-103), mirroring the
#     existing per-date block structure (15 action rows) already present
#     for every prior date (e.g. rows 74-88 for date 44542). ---

$labels = @(
    "total applicants",
    "withdraw before acceptance (WB)",
    "rejected groups",
    "preliminary rejection",
    "passive withdrawal",
    "rejection",
    "defer to MD app",
    "at least 1 MD/PhD acceptance",
    "available active",
    "request secondary",
    "interview scheduled",
    "available passive",
    "no action",
    "hold",
    "available"
)

# Only a few of the new rows carry an actual recorded count in column C;
# the rest (like most of the historical rows) are "#N/A" via =NA().
$counts = @{ 0 = 1765; 6 = 1; 7 = 329 }

$startRow = 89
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $startRow + $i

    # Copy the date cell's formatting (date number format, style) from the
    # last existing row so the new cell reuses the same style index instead
    # of minting a new one, then overwrite the value.
    $ws.Range("A88").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = 44555

    $ws.Range("B$r").Value = $labels[$i]

    if ($counts.ContainsKey($i)) {
        $ws.Range("C$r").Value = $counts[$i]
    } else {
        $ws.Range("C$r").Formula = "=NA()"
    }

    $ws.Range("D$r").Formula = "=NA()"
    $ws.Range("E$r").Formula = "=NA()"
}

# --- View-state: the user ended the session with the "Cycle_2021-2022"
#     sheet (3rd tab) selected/active, scrolled near the bottom of the new
#     data, instead of the "Cycle_2020-2021" sheet that was active before. ---
[void]$ws.Activate()
[void]$ws.Range("J89").Select()

Write-Output "applied fencer cycle updates"
